$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: new sampling week inserted at the top of the Apio series (rows 148-186),
# pushing the existing rows down by one; the last existing row becomes new row 187.

# Row 148
$ws.Range("A148").Value = 7
$ws.Range("B148").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C148").Value = "Ñuble"
$ws.Range("D148").Value = 44642
$ws.Range("E148").Value = 16
$ws.Range("F148").Value = 100112017
$ws.Range("G148").Value = "Apio"
$ws.Range("H148").Value = "Americana (o)"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 80
$ws.Range("K148").Value = 8000
$ws.Range("L148").Value = 9000
$ws.Range("M148").Value = 8500
$ws.Range("N148").Value = "`$/docena de matas"
$ws.Range("O148").Value = "Provincia del Elquí"
$ws.Range("P148").Value = 1417
$ws.Range("Q148").Value = 6
$ws.Range("R148").Value = "Hortaliza"

# Row 149
$ws.Range("A149").Value = 7
$ws.Range("B149").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C149").Value = "Ñuble"
$ws.Range("D149").Value = 44473
$ws.Range("E149").Value = 16
$ws.Range("F149").Value = 100112017
$ws.Range("G149").Value = "Apio"
$ws.Range("H149").Value = "Americana (o)"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 120
$ws.Range("K149").Value = 8000
$ws.Range("L149").Value = 9000
$ws.Range("M149").Value = 8500
$ws.Range("N149").Value = "`$/docena de matas"
$ws.Range("O149").Value = "Provincia del Elquí"
$ws.Range("P149").Value = 1417
$ws.Range("Q149").Value = 6
$ws.Range("R149").Value = "Hortaliza"

# Row 150
$ws.Range("A150").Value = 7
$ws.Range("B150").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C150").Value = "Ñuble"
$ws.Range("D150").Value = 44455
$ws.Range("E150").Value = 16
$ws.Range("F150").Value = 100112017
$ws.Range("G150").Value = "Apio"
$ws.Range("H150").Value = "Americana (o)"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 300
$ws.Range("K150").Value = 8500
$ws.Range("L150").Value = 9000
$ws.Range("M150").Value = 8750
$ws.Range("N150").Value = "`$/docena de matas"
$ws.Range("O150").Value = "Provincia del Elquí"
$ws.Range("P150").Value = 1458
$ws.Range("Q150").Value = 6
$ws.Range("R150").Value = "Hortaliza"

# Row 151
$ws.Range("A151").Value = 7
$ws.Range("B151").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C151").Value = "Ñuble"
$ws.Range("D151").Value = 44490
$ws.Range("E151").Value = 16
$ws.Range("F151").Value = 100112017
$ws.Range("G151").Value = "Apio"
$ws.Range("H151").Value = "Americana (o)"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 120
$ws.Range("K151").Value = 8000
$ws.Range("L151").Value = 8500
$ws.Range("M151").Value = 8250
$ws.Range("N151").Value = "`$/docena de matas"
$ws.Range("O151").Value = "Provincia del Elquí"
$ws.Range("P151").Value = 1375
$ws.Range("Q151").Value = 6
$ws.Range("R151").Value = "Hortaliza"

# Row 152
$ws.Range("A152").Value = 7
$ws.Range("B152").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C152").Value = "Ñuble"
$ws.Range("D152").Value = 44446
$ws.Range("E152").Value = 16
$ws.Range("F152").Value = 100112017
$ws.Range("G152").Value = "Apio"
$ws.Range("H152").Value = "Americana (o)"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 160
$ws.Range("K152").Value = 8000
$ws.Range("L152").Value = 8500
$ws.Range("M152").Value = 8250
$ws.Range("N152").Value = "`$/docena de matas"
$ws.Range("O152").Value = "Provincia del Elquí"
$ws.Range("P152").Value = 1375
$ws.Range("Q152").Value = 6
$ws.Range("R152").Value = "Hortaliza"

# Row 153
$ws.Range("A153").Value = 7
$ws.Range("B153").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C153").Value = "Ñuble"
$ws.Range("D153").Value = 44329
$ws.Range("E153").Value = 16
$ws.Range("F153").Value = 100112017
$ws.Range("G153").Value = "Apio"
$ws.Range("H153").Value = "Americana (o)"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 120
$ws.Range("K153").Value = 8000
$ws.Range("L153").Value = 9000
$ws.Range("M153").Value = 8500
$ws.Range("N153").Value = "`$/docena de matas"
$ws.Range("O153").Value = "Región de Coquimbo"
$ws.Range("P153").Value = 1417
$ws.Range("Q153").Value = 6
$ws.Range("R153").Value = "Hortaliza"

# Row 154
$ws.Range("A154").Value = 7
$ws.Range("B154").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C154").Value = "Ñuble"
$ws.Range("D154").Value = 44637
$ws.Range("E154").Value = 16
$ws.Range("F154").Value = 100112017
$ws.Range("G154").Value = "Apio"
$ws.Range("H154").Value = "Americana (o)"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 60
$ws.Range("K154").Value = 8000
$ws.Range("L154").Value = 8500
$ws.Range("M154").Value = 8250
$ws.Range("N154").Value = "`$/docena de matas"
$ws.Range("O154").Value = "Provincia del Elquí"
$ws.Range("P154").Value = 1375
$ws.Range("Q154").Value = 6
$ws.Range("R154").Value = "Hortaliza"

# Row 155
$ws.Range("A155").Value = 7
$ws.Range("B155").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C155").Value = "Ñuble"
$ws.Range("D155").Value = 44208
$ws.Range("E155").Value = 16
$ws.Range("F155").Value = 100112017
$ws.Range("G155").Value = "Apio"
$ws.Range("H155").Value = "Americana (o)"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 70
$ws.Range("K155").Value = 9500
$ws.Range("L155").Value = 10000
$ws.Range("M155").Value = 9821
$ws.Range("N155").Value = "`$/docena de matas"
$ws.Range("O155").Value = "Región de Coquimbo"
$ws.Range("P155").Value = 1637
$ws.Range("Q155").Value = 6
$ws.Range("R155").Value = "Hortaliza"

# Row 156
$ws.Range("A156").Value = 7
$ws.Range("B156").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C156").Value = "Ñuble"
$ws.Range("D156").Value = 44355
$ws.Range("E156").Value = 16
$ws.Range("F156").Value = 100112017
$ws.Range("G156").Value = "Apio"
$ws.Range("H156").Value = "Americana (o)"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 120
$ws.Range("K156").Value = 8000
$ws.Range("L156").Value = 8500
$ws.Range("M156").Value = 8250
$ws.Range("N156").Value = "`$/docena de matas"
$ws.Range("O156").Value = "Provincia del Elquí"
$ws.Range("P156").Value = 1375
$ws.Range("Q156").Value = 6
$ws.Range("R156").Value = "Hortaliza"

# Row 157
$ws.Range("A157").Value = 7
$ws.Range("B157").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C157").Value = "Ñuble"
$ws.Range("D157").Value = 44530
$ws.Range("E157").Value = 16
$ws.Range("F157").Value = 100112017
$ws.Range("G157").Value = "Apio"
$ws.Range("H157").Value = "Americana (o)"
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 100
$ws.Range("K157").Value = 8000
$ws.Range("L157").Value = 8500
$ws.Range("M157").Value = 8250
$ws.Range("N157").Value = "`$/docena de matas"
$ws.Range("O157").Value = "Provincia del Elquí"
$ws.Range("P157").Value = 1375
$ws.Range("Q157").Value = 6
$ws.Range("R157").Value = "Hortaliza"

# Row 158
$ws.Range("A158").Value = 7
$ws.Range("B158").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C158").Value = "Ñuble"
$ws.Range("D158").Value = 44483
$ws.Range("E158").Value = 16
$ws.Range("F158").Value = 100112017
$ws.Range("G158").Value = "Apio"
$ws.Range("H158").Value = "Americana (o)"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 160
$ws.Range("K158").Value = 8000
$ws.Range("L158").Value = 8500
$ws.Range("M158").Value = 8250
$ws.Range("N158").Value = "`$/docena de matas"
$ws.Range("O158").Value = "Provincia del Elquí"
$ws.Range("P158").Value = 1375
$ws.Range("Q158").Value = 6
$ws.Range("R158").Value = "Hortaliza"

# Row 159
$ws.Range("A159").Value = 7
$ws.Range("B159").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C159").Value = "Ñuble"
$ws.Range("D159").Value = 44617
$ws.Range("E159").Value = 16
$ws.Range("F159").Value = 100112017
$ws.Range("G159").Value = "Apio"
$ws.Range("H159").Value = "Americana (o)"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 60
$ws.Range("K159").Value = 8000
$ws.Range("L159").Value = 8500
$ws.Range("M159").Value = 8250
$ws.Range("N159").Value = "`$/docena de matas"
$ws.Range("O159").Value = "Provincia del Elquí"
$ws.Range("P159").Value = 1375
$ws.Range("Q159").Value = 6
$ws.Range("R159").Value = "Hortaliza"

# Row 160
$ws.Range("A160").Value = 7
$ws.Range("B160").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C160").Value = "Ñuble"
$ws.Range("D160").Value = 44489
$ws.Range("E160").Value = 16
$ws.Range("F160").Value = 100112017
$ws.Range("G160").Value = "Apio"
$ws.Range("H160").Value = "Americana (o)"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 180
$ws.Range("K160").Value = 8000
$ws.Range("L160").Value = 8500
$ws.Range("M160").Value = 8250
$ws.Range("N160").Value = "`$/docena de matas"
$ws.Range("O160").Value = "Provincia del Elquí"
$ws.Range("P160").Value = 1375
$ws.Range("Q160").Value = 6
$ws.Range("R160").Value = "Hortaliza"

# Row 161
$ws.Range("A161").Value = 7
$ws.Range("B161").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C161").Value = "Ñuble"
$ws.Range("D161").Value = 44264
$ws.Range("E161").Value = 16
$ws.Range("F161").Value = 100112017
$ws.Range("G161").Value = "Apio"
$ws.Range("H161").Value = "Americana (o)"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 65
$ws.Range("K161").Value = 8000
$ws.Range("L161").Value = 8500
$ws.Range("M161").Value = 8269
$ws.Range("N161").Value = "`$/docena de matas"
$ws.Range("O161").Value = "Región de Coquimbo"
$ws.Range("P161").Value = 1378
$ws.Range("Q161").Value = 6
$ws.Range("R161").Value = "Hortaliza"

# Row 162
$ws.Range("A162").Value = 7
$ws.Range("B162").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C162").Value = "Ñuble"
$ws.Range("D162").Value = 44396
$ws.Range("E162").Value = 16
$ws.Range("F162").Value = 100112017
$ws.Range("G162").Value = "Apio"
$ws.Range("H162").Value = "Americana (o)"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 120
$ws.Range("K162").Value = 8500
$ws.Range("L162").Value = 9000
$ws.Range("M162").Value = 8750
$ws.Range("N162").Value = "`$/docena de matas"
$ws.Range("O162").Value = "Provincia del Elquí"
$ws.Range("P162").Value = 1458
$ws.Range("Q162").Value = 6
$ws.Range("R162").Value = "Hortaliza"

# Row 163
$ws.Range("A163").Value = 7
$ws.Range("B163").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C163").Value = "Ñuble"
$ws.Range("D163").Value = 44279
$ws.Range("E163").Value = 16
$ws.Range("F163").Value = 100112017
$ws.Range("G163").Value = "Apio"
$ws.Range("H163").Value = "Americana (o)"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 80
$ws.Range("K163").Value = 8000
$ws.Range("L163").Value = 8000
$ws.Range("M163").Value = 8000
$ws.Range("N163").Value = "`$/docena de matas"
$ws.Range("O163").Value = "Región de Coquimbo"
$ws.Range("P163").Value = 1333
$ws.Range("Q163").Value = 6
$ws.Range("R163").Value = "Hortaliza"

# Row 164
$ws.Range("A164").Value = 7
$ws.Range("B164").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C164").Value = "Ñuble"
$ws.Range("D164").Value = 44330
$ws.Range("E164").Value = 16
$ws.Range("F164").Value = 100112017
$ws.Range("G164").Value = "Apio"
$ws.Range("H164").Value = "Americana (o)"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 120
$ws.Range("K164").Value = 7500
$ws.Range("L164").Value = 8000
$ws.Range("M164").Value = 7750
$ws.Range("N164").Value = "`$/docena de matas"
$ws.Range("O164").Value = "Región de Coquimbo"
$ws.Range("P164").Value = 1292
$ws.Range("Q164").Value = 6
$ws.Range("R164").Value = "Hortaliza"

# Row 165
$ws.Range("A165").Value = 7
$ws.Range("B165").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C165").Value = "Ñuble"
$ws.Range("D165").Value = 44504
$ws.Range("E165").Value = 16
$ws.Range("F165").Value = 100112017
$ws.Range("G165").Value = "Apio"
$ws.Range("H165").Value = "Americana (o)"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 120
$ws.Range("K165").Value = 8000
$ws.Range("L165").Value = 9000
$ws.Range("M165").Value = 8500
$ws.Range("N165").Value = "`$/docena de matas"
$ws.Range("O165").Value = "Provincia del Elquí"
$ws.Range("P165").Value = 1417
$ws.Range("Q165").Value = 6
$ws.Range("R165").Value = "Hortaliza"

# Row 166
$ws.Range("A166").Value = 7
$ws.Range("B166").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C166").Value = "Ñuble"
$ws.Range("D166").Value = 44257
$ws.Range("E166").Value = 16
$ws.Range("F166").Value = 100112017
$ws.Range("G166").Value = "Apio"
$ws.Range("H166").Value = "Americana (o)"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 55
$ws.Range("K166").Value = 7500
$ws.Range("L166").Value = 8000
$ws.Range("M166").Value = 7773
$ws.Range("N166").Value = "`$/docena de matas"
$ws.Range("O166").Value = "Región de Coquimbo"
$ws.Range("P166").Value = 1296
$ws.Range("Q166").Value = 6
$ws.Range("R166").Value = "Hortaliza"

# Row 167
$ws.Range("A167").Value = 7
$ws.Range("B167").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C167").Value = "Ñuble"
$ws.Range("D167").Value = 44301
$ws.Range("E167").Value = 16
$ws.Range("F167").Value = 100112017
$ws.Range("G167").Value = "Apio"
$ws.Range("H167").Value = "Americana (o)"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 120
$ws.Range("K167").Value = 8000
$ws.Range("L167").Value = 9000
$ws.Range("M167").Value = 8500
$ws.Range("N167").Value = "`$/docena de matas"
$ws.Range("O167").Value = "Región de Coquimbo"
$ws.Range("P167").Value = 1417
$ws.Range("Q167").Value = 6
$ws.Range("R167").Value = "Hortaliza"

# Row 168
$ws.Range("A168").Value = 7
$ws.Range("B168").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C168").Value = "Ñuble"
$ws.Range("D168").Value = 44370
$ws.Range("E168").Value = 16
$ws.Range("F168").Value = 100112017
$ws.Range("G168").Value = "Apio"
$ws.Range("H168").Value = "Americana (o)"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 120
$ws.Range("K168").Value = 7500
$ws.Range("L168").Value = 8000
$ws.Range("M168").Value = 7750
$ws.Range("N168").Value = "`$/docena de matas"
$ws.Range("O168").Value = "Provincia del Elquí"
$ws.Range("P168").Value = 1292
$ws.Range("Q168").Value = 6
$ws.Range("R168").Value = "Hortaliza"

# Row 169
$ws.Range("A169").Value = 7
$ws.Range("B169").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C169").Value = "Ñuble"
$ws.Range("D169").Value = 44487
$ws.Range("E169").Value = 16
$ws.Range("F169").Value = 100112017
$ws.Range("G169").Value = "Apio"
$ws.Range("H169").Value = "Americana (o)"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 180
$ws.Range("K169").Value = 8000
$ws.Range("L169").Value = 8500
$ws.Range("M169").Value = 8250
$ws.Range("N169").Value = "`$/docena de matas"
$ws.Range("O169").Value = "Provincia del Elquí"
$ws.Range("P169").Value = 1375
$ws.Range("Q169").Value = 6
$ws.Range("R169").Value = "Hortaliza"

# Row 170
$ws.Range("A170").Value = 7
$ws.Range("B170").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C170").Value = "Ñuble"
$ws.Range("D170").Value = 44385
$ws.Range("E170").Value = 16
$ws.Range("F170").Value = 100112017
$ws.Range("G170").Value = "Apio"
$ws.Range("H170").Value = "Americana (o)"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 120
$ws.Range("K170").Value = 8000
$ws.Range("L170").Value = 9000
$ws.Range("M170").Value = 8500
$ws.Range("N170").Value = "`$/docena de matas"
$ws.Range("O170").Value = "Provincia del Elquí"
$ws.Range("P170").Value = 1417
$ws.Range("Q170").Value = 6
$ws.Range("R170").Value = "Hortaliza"

# Row 171
$ws.Range("A171").Value = 7
$ws.Range("B171").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C171").Value = "Ñuble"
$ws.Range("D171").Value = 44413
$ws.Range("E171").Value = 16
$ws.Range("F171").Value = 100112017
$ws.Range("G171").Value = "Apio"
$ws.Range("H171").Value = "Americana (o)"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 120
$ws.Range("K171").Value = 8500
$ws.Range("L171").Value = 9000
$ws.Range("M171").Value = 8750
$ws.Range("N171").Value = "`$/docena de matas"
$ws.Range("O171").Value = "Provincia del Elquí"
$ws.Range("P171").Value = 1458
$ws.Range("Q171").Value = 6
$ws.Range("R171").Value = "Hortaliza"

# Row 172
$ws.Range("A172").Value = 7
$ws.Range("B172").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C172").Value = "Ñuble"
$ws.Range("D172").Value = 44272
$ws.Range("E172").Value = 16
$ws.Range("F172").Value = 100112017
$ws.Range("G172").Value = "Apio"
$ws.Range("H172").Value = "Americana (o)"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 80
$ws.Range("K172").Value = 8000
$ws.Range("L172").Value = 9000
$ws.Range("M172").Value = 8375
$ws.Range("N172").Value = "`$/docena de matas"
$ws.Range("O172").Value = "Región de Coquimbo"
$ws.Range("P172").Value = 1396
$ws.Range("Q172").Value = 6
$ws.Range("R172").Value = "Hortaliza"

# Row 173
$ws.Range("A173").Value = 7
$ws.Range("B173").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C173").Value = "Ñuble"
$ws.Range("D173").Value = 44299
$ws.Range("E173").Value = 16
$ws.Range("F173").Value = 100112017
$ws.Range("G173").Value = "Apio"
$ws.Range("H173").Value = "Americana (o)"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 120
$ws.Range("K173").Value = 8500
$ws.Range("L173").Value = 9000
$ws.Range("M173").Value = 8750
$ws.Range("N173").Value = "`$/docena de matas"
$ws.Range("O173").Value = "Región de Coquimbo"
$ws.Range("P173").Value = 1458
$ws.Range("Q173").Value = 6
$ws.Range("R173").Value = "Hortaliza"

# Row 174
$ws.Range("A174").Value = 7
$ws.Range("B174").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C174").Value = "Ñuble"
$ws.Range("D174").Value = 44610
$ws.Range("E174").Value = 16
$ws.Range("F174").Value = 100112017
$ws.Range("G174").Value = "Apio"
$ws.Range("H174").Value = "Americana (o)"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 60
$ws.Range("K174").Value = 8000
$ws.Range("L174").Value = 8500
$ws.Range("M174").Value = 8250
$ws.Range("N174").Value = "`$/docena de matas"
$ws.Range("O174").Value = "Provincia del Elquí"
$ws.Range("P174").Value = 1375
$ws.Range("Q174").Value = 6
$ws.Range("R174").Value = "Hortaliza"

# Row 175
$ws.Range("A175").Value = 7
$ws.Range("B175").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C175").Value = "Ñuble"
$ws.Range("D175").Value = 44312
$ws.Range("E175").Value = 16
$ws.Range("F175").Value = 100112017
$ws.Range("G175").Value = "Apio"
$ws.Range("H175").Value = "Americana (o)"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 120
$ws.Range("K175").Value = 8000
$ws.Range("L175").Value = 9000
$ws.Range("M175").Value = 8500
$ws.Range("N175").Value = "`$/docena de matas"
$ws.Range("O175").Value = "Provincia del Elquí"
$ws.Range("P175").Value = 1417
$ws.Range("Q175").Value = 6
$ws.Range("R175").Value = "Hortaliza"

# Row 176
$ws.Range("A176").Value = 7
$ws.Range("B176").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C176").Value = "Ñuble"
$ws.Range("D176").Value = 44399
$ws.Range("E176").Value = 16
$ws.Range("F176").Value = 100112017
$ws.Range("G176").Value = "Apio"
$ws.Range("H176").Value = "Americana (o)"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 120
$ws.Range("K176").Value = 8500
$ws.Range("L176").Value = 9000
$ws.Range("M176").Value = 8750
$ws.Range("N176").Value = "`$/docena de matas"
$ws.Range("O176").Value = "Provincia del Elquí"
$ws.Range("P176").Value = 1458
$ws.Range("Q176").Value = 6
$ws.Range("R176").Value = "Hortaliza"

# Row 177
$ws.Range("A177").Value = 7
$ws.Range("B177").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C177").Value = "Ñuble"
$ws.Range("D177").Value = 44522
$ws.Range("E177").Value = 16
$ws.Range("F177").Value = 100112017
$ws.Range("G177").Value = "Apio"
$ws.Range("H177").Value = "Americana (o)"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 80
$ws.Range("K177").Value = 8000
$ws.Range("L177").Value = 8500
$ws.Range("M177").Value = 8250
$ws.Range("N177").Value = "`$/docena de matas"
$ws.Range("O177").Value = "Provincia del Elquí"
$ws.Range("P177").Value = 1375
$ws.Range("Q177").Value = 6
$ws.Range("R177").Value = "Hortaliza"

# Row 178
$ws.Range("A178").Value = 7
$ws.Range("B178").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C178").Value = "Ñuble"
$ws.Range("D178").Value = 44543
$ws.Range("E178").Value = 16
$ws.Range("F178").Value = 100112017
$ws.Range("G178").Value = "Apio"
$ws.Range("H178").Value = "Americana (o)"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 100
$ws.Range("K178").Value = 8000
$ws.Range("L178").Value = 8500
$ws.Range("M178").Value = 8250
$ws.Range("N178").Value = "`$/docena de matas"
$ws.Range("O178").Value = "Provincia del Elquí"
$ws.Range("P178").Value = 1375
$ws.Range("Q178").Value = 6
$ws.Range("R178").Value = "Hortaliza"

# Row 179
$ws.Range("A179").Value = 7
$ws.Range("B179").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C179").Value = "Ñuble"
$ws.Range("D179").Value = 44277
$ws.Range("E179").Value = 16
$ws.Range("F179").Value = 100112017
$ws.Range("G179").Value = "Apio"
$ws.Range("H179").Value = "Americana (o)"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 120
$ws.Range("K179").Value = 8000
$ws.Range("L179").Value = 9000
$ws.Range("M179").Value = 8500
$ws.Range("N179").Value = "`$/docena de matas"
$ws.Range("O179").Value = "Región de Coquimbo"
$ws.Range("P179").Value = 1417
$ws.Range("Q179").Value = 6
$ws.Range("R179").Value = "Hortaliza"

# Row 180
$ws.Range("A180").Value = 7
$ws.Range("B180").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C180").Value = "Ñuble"
$ws.Range("D180").Value = 44258
$ws.Range("E180").Value = 16
$ws.Range("F180").Value = 100112017
$ws.Range("G180").Value = "Apio"
$ws.Range("H180").Value = "Americana (o)"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 100
$ws.Range("K180").Value = 7000
$ws.Range("L180").Value = 8000
$ws.Range("M180").Value = 7600
$ws.Range("N180").Value = "`$/docena de matas"
$ws.Range("O180").Value = "Región de Coquimbo"
$ws.Range("P180").Value = 1267
$ws.Range("Q180").Value = 6
$ws.Range("R180").Value = "Hortaliza"

# Row 181
$ws.Range("A181").Value = 7
$ws.Range("B181").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C181").Value = "Ñuble"
$ws.Range("D181").Value = 44390
$ws.Range("E181").Value = 16
$ws.Range("F181").Value = 100112017
$ws.Range("G181").Value = "Apio"
$ws.Range("H181").Value = "Americana (o)"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 120
$ws.Range("K181").Value = 8000
$ws.Range("L181").Value = 9000
$ws.Range("M181").Value = 8500
$ws.Range("N181").Value = "`$/docena de matas"
$ws.Range("O181").Value = "Provincia del Elquí"
$ws.Range("P181").Value = 1417
$ws.Range("Q181").Value = 6
$ws.Range("R181").Value = "Hortaliza"

# Row 182
$ws.Range("A182").Value = 7
$ws.Range("B182").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C182").Value = "Ñuble"
$ws.Range("D182").Value = 44349
$ws.Range("E182").Value = 16
$ws.Range("F182").Value = 100112017
$ws.Range("G182").Value = "Apio"
$ws.Range("H182").Value = "Americana (o)"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 160
$ws.Range("K182").Value = 7500
$ws.Range("L182").Value = 8000
$ws.Range("M182").Value = 7750
$ws.Range("N182").Value = "`$/docena de matas"
$ws.Range("O182").Value = "Región de Coquimbo"
$ws.Range("P182").Value = 1292
$ws.Range("Q182").Value = 6
$ws.Range("R182").Value = "Hortaliza"

# Row 183
$ws.Range("A183").Value = 7
$ws.Range("B183").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C183").Value = "Ñuble"
$ws.Range("D183").Value = 44498
$ws.Range("E183").Value = 16
$ws.Range("F183").Value = 100112017
$ws.Range("G183").Value = "Apio"
$ws.Range("H183").Value = "Americana (o)"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 120
$ws.Range("K183").Value = 8000
$ws.Range("L183").Value = 9000
$ws.Range("M183").Value = 8500
$ws.Range("N183").Value = "`$/docena de matas"
$ws.Range("O183").Value = "Provincia del Elquí"
$ws.Range("P183").Value = 1417
$ws.Range("Q183").Value = 6
$ws.Range("R183").Value = "Hortaliza"

# Row 184
$ws.Range("A184").Value = 7
$ws.Range("B184").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C184").Value = "Ñuble"
$ws.Range("D184").Value = 44418
$ws.Range("E184").Value = 16
$ws.Range("F184").Value = 100112017
$ws.Range("G184").Value = "Apio"
$ws.Range("H184").Value = "Americana (o)"
$ws.Range("I184").Value = "Primera"
$ws.Range("J184").Value = 120
$ws.Range("K184").Value = 8500
$ws.Range("L184").Value = 9000
$ws.Range("M184").Value = 8750
$ws.Range("N184").Value = "`$/docena de matas"
$ws.Range("O184").Value = "Provincia del Elquí"
$ws.Range("P184").Value = 1458
$ws.Range("Q184").Value = 6
$ws.Range("R184").Value = "Hortaliza"

# Row 185
$ws.Range("A185").Value = 7
$ws.Range("B185").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C185").Value = "Ñuble"
$ws.Range("D185").Value = 44628
$ws.Range("E185").Value = 16
$ws.Range("F185").Value = 100112017
$ws.Range("G185").Value = "Apio"
$ws.Range("H185").Value = "Americana (o)"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 60
$ws.Range("K185").Value = 8500
$ws.Range("L185").Value = 9000
$ws.Range("M185").Value = 8750
$ws.Range("N185").Value = "`$/docena de matas"
$ws.Range("O185").Value = "Provincia del Elquí"
$ws.Range("P185").Value = 1458
$ws.Range("Q185").Value = 6
$ws.Range("R185").Value = "Hortaliza"

# Row 186
$ws.Range("A186").Value = 7
$ws.Range("B186").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C186").Value = "Ñuble"
$ws.Range("D186").Value = 44335
$ws.Range("E186").Value = 16
$ws.Range("F186").Value = 100112017
$ws.Range("G186").Value = "Apio"
$ws.Range("H186").Value = "Americana (o)"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 120
$ws.Range("K186").Value = 8000
$ws.Range("L186").Value = 8500
$ws.Range("M186").Value = 8250
$ws.Range("N186").Value = "`$/docena de matas"
$ws.Range("O186").Value = "Región de Coquimbo"
$ws.Range("P186").Value = 1375
$ws.Range("Q186").Value = 6
$ws.Range("R186").Value = "Hortaliza"

# Row 187
$ws.Range("A187").Value = 7
$ws.Range("B187").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C187").Value = "Ñuble"
$ws.Range("D187").Value = 44552
$ws.Range("E187").Value = 16
$ws.Range("F187").Value = 100112017
$ws.Range("G187").Value = "Apio"
$ws.Range("H187").Value = "Americana (o)"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 60
$ws.Range("K187").Value = 8000
$ws.Range("L187").Value = 8500
$ws.Range("M187").Value = 8250
$ws.Range("N187").Value = "`$/docena de matas"
$ws.Range("O187").Value = "Provincia del Elquí"
$ws.Range("P187").Value = 1375
$ws.Range("Q187").Value = 6
$ws.Range("R187").Value = "Hortaliza"

$ws.Range("A1").Select()
